# [IMPL] CRUD FACTURAS REALIZADO
# - Modelo de facturas implementado
# - Mostrar facturas (CLI), listar facturas (ADMIN), mostrar facturas (ADMIN)
#
# This script reproduces, via Excel COM-interop, the edits made to
# "Product-Backlog.xlsx": a handful of SPRINT-BACKLOG task rows (53-59) get
# their "Fecha Fin Real" (column D) filled in and their "Estado" (column F)
# flipped to REALIZADO, two brand-new tasks about facturas views are added,
# and the view/selection state of the two sheets is updated to reflect where
# the author was last working.

$wb = $excel.ActiveWorkbook
$wsSprint = $wb.Worksheets.Item("SPRINT-BACKLOG")
$wsAux = $wb.Worksheets.Item("AUX")

# ---------------------------------------------------------------------------
# 1. SPRINT-BACKLOG data edits (rows 53-59)
# ---------------------------------------------------------------------------

# Row 53: "Creación de la vista de árbol  de facturas" task finished
$wsSprint.Range("D53").Value = 44263
$wsSprint.Range("F53").Value = "REALIZADO"

# Row 54: new task "Creación de la vista de árbol  de facturas"
$wsSprint.Range("A54").Value = "Creación de la vista de árbol  de facturas"
$wsSprint.Range("D54").Value = 44264
$wsSprint.Range("F54").Value = "REALIZADO"

# Row 55: new task "Creación de la vista de formulario de facturas"
$wsSprint.Range("A55").Value = "Creación de la vista de formulario de facturas"
$wsSprint.Range("D55").Value = 44264
$wsSprint.Range("F55").Value = "REALIZADO"

# Row 56
$wsSprint.Range("D56").Value = 44264
$wsSprint.Range("F56").Value = "REALIZADO"

# Row 57
$wsSprint.Range("D57").Value = 44264
$wsSprint.Range("F57").Value = "REALIZADO"

# Row 58 (status left as-is, only the real end date gets filled in)
$wsSprint.Range("D58").Value = 44264

# Row 59
$wsSprint.Range("D59").Value = 44264
$wsSprint.Range("F59").Value = "REALIZADO"

# Recalculate so the AUX sheet's SUMIF-driven burn-down figures refresh.
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------------
# 2. Column width tweak on SPRINT-BACKLOG (new explicit width for column F)
# ---------------------------------------------------------------------------
$wsSprint.Columns.Item(6).ColumnWidth = 11.7

# ---------------------------------------------------------------------------
# 3. View / selection state
# ---------------------------------------------------------------------------

# AUX: selection moved to J12, scrolled back to the top-left corner.
$wsAux.Activate()
$wsAux.Range("J12").Select()

# SPRINT-BACKLOG: back on top as the active sheet, zoomed to 85%, selection
# on D60.
$wsSprint.Activate()
$wsSprint.Range("D60").Select()
$excel.ActiveWindow.Zoom = 85
